$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
$text = @'
88
B C E I n c . 
 2 0 1 5 A N N U A L R E PO R T
7 Selected annual 
and quarterly information
7.1 Annual financial information
The following table shows selected consolidated financial data of BCE for 2015, 2014 and 2013, prepared in accordance with IFRS as issued 
by the International Accounting Standards Board (IASB). We discuss the factors that caused our results to vary over the past two years 
throughout this MD&A.
2015
2014 (1)
2013 (2)
CONSOLIDATED INCOME STATEMENTS
Operating revenues
21,514
21,042
20,400
Operating costs
(12,963)
(12,739)
(12,311)
Adjusted EBITDA
8,551
8,303
8,089
Severance, acquisition and other costs
(446)
(216)
(406)
Depreciation
(2,890)
(2,880)
(2,734)
Amortization
(530)
(572)
(646)
Finance costs
Interest expense
(909)
(929)
(931)
Interest on post-employment benefit obligations
(110)
(101)
(150)
Other (expense) income
(12)
42
(6)
Income taxes
(924)
(929)
(828)
Net earnings
2,730
2,718
2,388
Net earnings attributable to:
Common shareholders
2,526
2,363
1,975
Preferred shareholders
152
137
131
Non-controlling interest
52
218
282
Net earnings
2,730
2,718
2,388
Net earnings per common share
Basic
2.98
2.98
2.55
Diluted
2.98
2.97
2.54
Included in net earnings:
Severance, acquisition and other costs
(327)
(148)
(299)
Net gains (losses) on investments
21
8
(7)
Early debt redemption costs
(13)
(21)
(36)
Adjusted net earnings
2,845
2,524
2,317
Adjusted EPS
3.36
3.18
2.99
RATIOS
Adjusted EBITDA margin (%)
39.7%
39.5%
39.7%
Return on equity (%) (3)
21.1%
21.0%
17.9%
(1) On October 31, 2014, BCE completed its acquisition of all the issued and outstanding common shares of Bell Aliant that it did not already own. Refer to section 6.5, Privatization 
of Bell Aliant for further details on the transaction.
(2) On July 5, 2013, BCE acquired 100% of the issued and outstanding shares of Astral. As part of its approval of the Astral acquisition, the CRTC ordered BCE to spend $246.9 million in 
new benefits for French- and English-language TV, radio and film content development, support for emerging Canadian musical talent, training and professional development for 
Canadian media, and new consumer participation initiatives. The present value of this tangible benefits obligation, amounting to $230 million, was recorded as an acquisition cost in 
Severance, acquisition and other costs in 2013. Total acquisition and other costs relating to Astral, including the tangible benefits obligation, amounted to $266 million in 2013.
(3) Net earnings attributable to common shareholders divided by total average equity attributable to BCE shareholders excluding preferred shares.

'@
$ws.Range("A2").Value = $text
$ws.Rows.Item(2).EntireRow.AutoFit()
